$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Fri Oct 25 11:30:29 EDT 2024",
    "Fri Oct 25 11:30:39 EDT 2024",
    "Fri Oct 25 11:30:49 EDT 2024",
    "Fri Oct 25 11:31:00 EDT 2024",
    "Fri Oct 25 11:31:11 EDT 2024",
    "Fri Oct 25 11:31:21 EDT 2024",
    "Fri Oct 25 11:31:31 EDT 2024",
    "Fri Oct 25 11:31:42 EDT 2024",
    "Fri Oct 25 11:31:52 EDT 2024",
    "Fri Oct 25 11:32:02 EDT 2024",
    "Fri Oct 25 11:32:13 EDT 2024",
    "Fri Oct 25 11:32:23 EDT 2024",
    "Fri Oct 25 11:32:33 EDT 2024",
    "Fri Oct 25 11:32:43 EDT 2024",
    "Fri Oct 25 11:32:53 EDT 2024",
    "Fri Oct 25 11:33:03 EDT 2024",
    "Fri Oct 25 11:33:14 EDT 2024",
    "Fri Oct 25 11:33:24 EDT 2024",
    "Fri Oct 25 11:33:34 EDT 2024",
    "Fri Oct 25 11:33:44 EDT 2024",
    "Fri Oct 25 11:33:54 EDT 2024",
    "Fri Oct 25 11:34:05 EDT 2024",
    "Fri Oct 25 11:34:15 EDT 2024",
    "Fri Oct 25 11:34:25 EDT 2024",
    "Fri Oct 25 11:34:36 EDT 2024",
    "Fri Oct 25 11:34:47 EDT 2024",
    "Fri Oct 25 11:34:57 EDT 2024",
    "Fri Oct 25 11:35:08 EDT 2024",
    "Fri Oct 25 11:35:18 EDT 2024"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
